$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("5002020-12.2020.8.21.0067", "9000539-43.2020.8.21.0067", "Migrado"),
    @("5002025-34.2020.8.21.0067", "9000434-66.2020.8.21.0067", "Migrado"),
    @("5000995-32.2018.8.21.0067", "9000796-39.2018.8.21.0067", "Migrado"),
    @("5002032-26.2020.8.21.0067", "9000450-20.2020.8.21.0067", "Migrado"),
    @("5000229-27.2011.8.21.0001", "0420539-74.2011.8.21.0001", "Digitalizado"),
    @("5000297-74.2011.8.21.0001", "0413928-08.2011.8.21.0001", "Digitalizado"),
    @("5000320-49.2013.8.21.0001", "0044684-94.2013.8.21.0001", "Digitalizado"),
    @("5003815-38.2012.8.21.0001", "0023925-46.2012.8.21.0001", "Digitalizado"),
    @("5002031-41.2020.8.21.0067", "9000380-03.2020.8.21.0067", "Migrado"),
    @("5002033-11.2020.8.21.0067", "9000365-34.2020.8.21.0067", "Migrado"),
    @("5002037-48.2020.8.21.0067", "9000624-29.2020.8.21.0067", "Migrado"),
    @("5002090-63.2019.8.21.0067", "9001254-22.2019.8.21.0067", "Migrado"),
    @("5002030-56.2020.8.21.0067", "9000405-16.2020.8.21.0067", "Migrado"),
    @("5002091-48.2019.8.21.0067", "9001270-73.2019.8.21.0067", "Migrado"),
    @("5092615-95.2019.8.21.0001", "9015530-62.2019.8.21.0001", "Migrado")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
